# "dashboard a filter working"
# Replace the old per-agent RITM list with a single fresh row (RITM0496748)
# and add three new "dashboard" columns: Dashboard Status, Present Time,
# Closed Time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop every existing hyperlink (D2:D8 each carried one) before we touch
# the data, otherwise their relationship entries linger.
$ws.Cells.Hyperlinks.Delete()

# Wipe out all of the old agent rows (2-8), but leave D2 alone for a moment
# so its "Hyperlink" cell style (already applied in the template) survives.
$ws.Range("A2:C8").Clear()
$ws.Range("E2:K8").Clear()
$ws.Range("D3:D8").Clear()

# New single data row: just the RITM number for the dashboard filter test.
$ws.Range("D2").Value = "RITM0496748"

# Re-attach a hyperlink (ServiceNow request item) to the new RITM number,
# then re-assert the Hyperlink cell style since adding the link resets it.
$ws.Hyperlinks.Add($ws.Range("D2"), "https://europarl.service-now.com/nav_to.do?uri=%2Fsc_req_item.do%3Fsys_id%3D30877432d1026706d7e805da846a32c3%26sysparm_view%3D")
$ws.Range("D2").Style = "Hyperlink"

# New dashboard columns appended after the existing headers.
$ws.Range("L1").Value = "Dashboard Status"
$ws.Range("M1").Value = "Present Time"
$ws.Range("N1").Value = "Closed Time"

# Leave the selection where the user ended up while testing the filter.
$ws.Range("E2").Select()
